{"js": "// Office.js (Word JavaScript API) edit script.\n// Body: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Remove the \"\u62a5\u544a\u65e5\u671f\uff1a...\" paragraph and the blank paragraph that\n//    immediately follows it (both sit right after the \"\u793e\u4ea4\u5a92\u4f53\u4e0a\u7684\u65b0\u70ed\u70b9\"\n//    heading, before the \"\u6700\u8fd1\u51e0\u5468...\" paragraph).\nconst dateParaText = \"\u62a5\u544a\u65e5\u671f\uff1a2024 \u5e74 1 \u6708 22 \u65e5\";\nlet dateParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === dateParaText) {\n    dateParaIndex = i;\n    break;\n  }\n}\n\nif (dateParaIndex !== -1) {\n  const blankParaIndex = dateParaIndex + 1;\n  // Delete the blank paragraph first so the date paragraph's index stays valid.\n  if (blankParaIndex < paragraphs.items.length && paragraphs.items[blankParaIndex].text === \" \") {\n    paragraphs.items[blankParaIndex].delete();\n  }\n  paragraphs.items[dateParaIndex].delete();\n  await context.sync();\n}\n\n// 2) Update the closing-paragraph sentence with the revised wording.\nconst oldSentence = \"\u5982\u679c\u8be5\u4ea7\u54c1\u80fd\u591f\u5728\u5065\u5eb7\u548c\u5065\u8eab\u9886\u57df\u7ee7\u7eed\u4fdd\u6301\u76ee\u524d\u7684\u9500\u552e\u901f\u5ea6\uff0c\u90a3\u4e48\u5b83\u5c31\u53ef\u80fd\u5df2\u7ecf\u51c6\u5907\u597d\u5728\u5168\u56fd\u8303\u56f4\u5185\u63a8\u51fa\u3002\";\nconst newSentence = \"\u5982\u679c\u8be5\u4ea7\u54c1\u80fd\u591f\u5728\u5065\u5eb7\u4e0e\u5065\u8eab\u9886\u57df\u7ee7\u7eed\u7ef4\u6301\u5f53\u524d\u7684\u9500\u552e\u52bf\u5934\uff0c\u90a3\u4e48\u5b83\u53ef\u80fd\u5df2\u7ecf\u51c6\u5907\u597d\u5728\u5168\u56fd\u8303\u56f4\u5185\u63a8\u5e7f\u3002\";\n\nconst searchResults = body.search(oldSentence, { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(newSentence, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the \"\u62a5\u544a\u65e5\u671f\uff1a...\" paragraph and the blank paragraph that\n#    immediately follows it (both sit right after the \"\u793e\u4ea4\u5a92\u4f53\u4e0a\u7684\u65b0\u70ed\u70b9\"\n#    heading, before the \"\u6700\u8fd1\u51e0\u5468...\" paragraph).\n$dateParaText = \"\u62a5\u544a\u65e5\u671f\uff1a2024 \u5e74 1 \u6708 22 \u65e5\"\n\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $dateParaText) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -ne -1) {\n    # Delete the blank paragraph that follows first (if present), then the\n    # date paragraph itself, so indices stay valid while deleting.\n    $nextIndex = $targetIndex + 1\n    if ($nextIndex -le $d.Paragraphs.Count) {\n        $nextText = $d.Paragraphs.Item($nextIndex).Range.Text.TrimEnd([char]13, [char]7)\n        if ($nextText -eq \" \") {\n            $d.Paragraphs.Item($nextIndex).Range.Delete()\n        }\n    }\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n}\n\n# 2) Update the closing-paragraph sentence with the revised wording.\n$oldSentence = \"\u5982\u679c\u8be5\u4ea7\u54c1\u80fd\u591f\u5728\u5065\u5eb7\u548c\u5065\u8eab\u9886\u57df\u7ee7\u7eed\u4fdd\u6301\u76ee\u524d\u7684\u9500\u552e\u901f\u5ea6\uff0c\u90a3\u4e48\u5b83\u5c31\u53ef\u80fd\u5df2\u7ecf\u51c6\u5907\u597d\u5728\u5168\u56fd\u8303\u56f4\u5185\u63a8\u51fa\u3002\"\n$newSentence = \"\u5982\u679c\u8be5\u4ea7\u54c1\u80fd\u591f\u5728\u5065\u5eb7\u4e0e\u5065\u8eab\u9886\u57df\u7ee7\u7eed\u7ef4\u6301\u5f53\u524d\u7684\u9500\u552e\u52bf\u5934\uff0c\u90a3\u4e48\u5b83\u53ef\u80fd\u5df2\u7ecf\u51c6\u5907\u597d\u5728\u5168\u56fd\u8303\u56f4\u5185\u63a8\u5e7f\u3002\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldSentence\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newSentence\n$find.Execute(\n    $oldSentence,  # FindText\n    $false,        # MatchCase\n    $false,        # MatchWholeWord\n    $false,        # MatchWildcards\n    $false,        # MatchSoundsLike\n    $false,        # MatchAllWordForms\n    $true,         # Forward\n    1,             # Wrap (wdFindContinue)\n    $false,        # Format\n    $newSentence,  # ReplaceWith\n    2              # Replace (wdReplaceAll)\n)\n"}
